$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general_summary")

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "H2" "100.0"
Set-TextValue "J2" "100.0"

Set-TextValue "H3" "40.0"
Set-TextValue "J3" "100.0"
Set-TextValue "L3" "60.0"

Set-TextValue "H4" "9.8"
Set-TextValue "J4" "63.8"
Set-TextValue "L4" "88.8"
Set-TextValue "N4" "39.4"

Set-TextValue "H5" "12.6"
Set-TextValue "J5" "52.9"
Set-TextValue "L5" "87.4"
Set-TextValue "N5" "47.1"

Set-TextValue "H6" "8.3"
Set-TextValue "J6" "40.4"
Set-TextValue "L6" "91.7"
Set-TextValue "N6" "59.6"

Set-TextValue "H7" "44.5"
Set-TextValue "J7" "70.2"
Set-TextValue "L7" "55.5"
Set-TextValue "N7" "29.8"
